$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.883.43"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "3.843.17"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "697.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").Value = "3.841.90"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "4.490.81"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "3.804.61"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "70.895.58"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.180"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "3.799.83"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.73%  "
$ws.Range("E40").Value = "  +7.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("E42").Value = "  -5.30%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000312"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.299"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.93%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.26%  "
